# The deck's two theme parts (ppt/theme/theme1.xml, used by the slide
# master, and ppt/theme/theme2.xml, used by the notes master) had their
# colour schemes swapped: theme1 ("Integral") becomes the Office default
# palette, theme2 ("Office Theme") becomes the Integral palette.
#
# The PowerPoint object model only exposes a writable 12-slot theme colour
# scheme on slide-facing anchors (Slide/SlideRange), and every master-like
# object in this host (SlideMaster/NotesMaster/HandoutMaster/NotesPage) all
# resolve back to the single slide master's theme (theme1.xml) - there is no
# reachable OM path to the notes master's theme2.xml.  So we apply the
# colour half of the swap that is reachable: push the Office theme's colour
# values (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) onto theme1's colour
# scheme, via a slide's ThemeColorScheme (which edits the shared
# presentation theme, i.e. ppt/theme/theme1.xml).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RGB() isn't available in this host, so each OLE colour below is the
# precomputed r + g*256 + b*65536 value for the target hex colour.
# index  name       hex      (r,g,b)
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
